{"js": "// Append a \"BUGS\" section with a struck-through bullet list (reusing the\n// existing numId=7 list used by \"IDEAS\") at the very end of the document\n// body, mirroring the target OOXML diff.\n\nconst body = context.document.body;\n\n// Office.js merges OOXML inserted at Word.InsertLocation.end of the body\n// into the last existing paragraph instead of creating new ones (there is\n// no paragraph mark after the last paragraph before the sectPr). To avoid\n// that, first append a fresh placeholder paragraph, then replace that\n// paragraph (and only that paragraph) with the full OOXML fragment below.\nconst placeholder = body.insertParagraph(\"\", Word.InsertLocation.end);\nawait context.sync();\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>\n</pkg:xmlData></pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>\n<w:p><w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>BUGS</w:t></w:r></w:p>\n<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"7\"/></w:numPr><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>Fix prices</w:t></w:r><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> (changed formula from N*N/2 to </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>catn</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>=cat(n-</w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>1)</w:t></w:r><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>*</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>2+</w:t></w:r><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>cat(n-1)</w:t></w:r><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>/1,5</w:t></w:r></w:p>\n<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"7\"/></w:numPr><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>Implement decimals</w:t></w:r></w:p>\n<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"7\"/></w:numPr><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>Fix fish</w:t></w:r></w:p>\n</w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>`;\n\nplaceholder.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Append a \"BUGS\" section with a struck-through bullet list (reusing the\n# existing numId=7 list used by \"IDEAS\") at the very end of the document,\n# mirroring the target OOXML diff.\n\n$d = $word.ActiveDocument\n\n$end = $d.Content\n$end.Collapse(0)  # wdCollapseEnd\n\n$ooxml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>\n</pkg:xmlData></pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>\n<w:p><w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>BUGS</w:t></w:r></w:p>\n<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"7\"/></w:numPr><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>Fix prices</w:t></w:r><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> (changed formula from N*N/2 to </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>catn</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>=cat(n-</w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>1)</w:t></w:r><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>*</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>2+</w:t></w:r><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>cat(n-1)</w:t></w:r><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>/1,5</w:t></w:r></w:p>\n<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"7\"/></w:numPr><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>Implement decimals</w:t></w:r></w:p>\n<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"7\"/></w:numPr><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:lang w:val=\"en-US\"/></w:rPr><w:t>Fix fish</w:t></w:r></w:p>\n</w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>\n'@\n\n$end.InsertXML($ooxml)\n"}
